$wb = $excel.ActiveWorkbook

# --- Update the time_taken (F column) timestamps on the "data" sheet ---
$ds = $wb.Worksheets.Item("data")
$ds.Range("F2").Value = "2021-10-05 14:19:21.521073"
$ds.Range("F3").Value = "2021-10-05 14:19:21.521082"
$ds.Range("F4").Value = "2021-10-05 14:19:21.521085"
$ds.Range("F5").Value = "2021-10-05 14:19:21.521088"
$ds.Range("F6").Value = "2021-10-05 14:19:21.521091"

# --- Add the new "metadata" sheet right after "data" ---
$ms = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ds)
$ms.Name = "metadata"

# Header row
$ms.Range("B1").Value = "data_name"
$ms.Range("C1").Value = "data_id"
$ms.Range("D1").Value = "data_version"
$ms.Range("E1").Value = "data_version_created"
$ms.Range("F1").Value = "panel_query_time"
$ms.Range("G1").Value = "panel_get_request"

# Header formatting: bold, thin border, centered horizontally, top vertically
$headerRange = $ms.Range("B1:G1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# Data row
$ms.Range("A2").Value = 0
$ms.Range("B2").Value = "Breast cancer pertinent cancer susceptibility"
$ms.Range("C2").Value = 55
$ms.Range("D2").Value = "'1.3"
$ms.Range("E2").Value = "2021-09-29T13:16:42.827685Z"
$ms.Range("F2").Value = "2021-10-05 14:19:21.517550"
$ms.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/55/?format=json"

# A2 (index column) uses the same bold/border/aligned style as the header
$a2 = $ms.Range("A2")
$a2.Font.Bold = $true
$a2.Borders.LineStyle = 1
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160

# Keep "data" as the active sheet (matches original workbook view state)
$ds.Activate()
